# The catalog's first two columns trade places: the "study" name (used
# to be column B) now leads as column A, and the numeric "id" (used to
# be column A) moves to column B. Cutting the whole column B and
# inserting it before column A performs the swap while carrying each
# column's own width/format along with it (so column A inherits column
# B's bestFit width, matching the wider "study" text it now holds).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$studyColumn = $ws.Columns.Item(2)
$studyColumn.Cut() | Out-Null
$ws.Columns.Item(1).Insert() | Out-Null

# Update the active selection to match the post-edit cursor position.
$ws.Range("J28").Select() | Out-Null
